$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 15).Value = 44.9833693541098
$ws.Cells.Item(4, 15).Value = 41.9605705904505
$ws.Cells.Item(5, 15).Value = 45.037318436426
$ws.Cells.Item(6, 15).Value = 44.4455922444536
$ws.Cells.Item(7, 15).Value = 10.3622362790511
$ws.Cells.Item(8, 15).Value = 39.8111704972711
$ws.Cells.Item(9, 15).Value = 33.9574450767427
$ws.Cells.Item(10, 15).Value = 31.1119906554538
$ws.Cells.Item(11, 15).Value = 52.9507022420666
$ws.Cells.Item(12, 15).Value = 18.26842010274
$ws.Cells.Item(13, 15).Value = 34.5024513087645
$ws.Cells.Item(14, 15).Value = 8.12233191914566
$ws.Cells.Item(15, 15).Value = 16.4418060108669
$ws.Cells.Item(16, 15).Value = 13.6921460620691
$ws.Cells.Item(17, 15).Value = 34.55980686207
$ws.Cells.Item(18, 15).Value = 66.0522012814459
$ws.Cells.Item(19, 15).Value = 38.9930717310614
$ws.Cells.Item(20, 15).Value = 44.6241623191673
$ws.Cells.Item(21, 15).Value = 39.8882587328187
$ws.Cells.Item(22, 15).Value = 15.3632638728612
$ws.Cells.Item(23, 15).Value = 32.6669339380576
$ws.Cells.Item(24, 15).Value = 11.5195131831656
$ws.Cells.Item(25, 15).Value = 149.333123346009
$ws.Cells.Item(26, 15).Value = ".."
$ws.Cells.Item(27, 15).Value = 8.5738621328595
$ws.Cells.Item(28, 15).Value = 10.7623285479165
$ws.Cells.Item(29, 15).Value = 23.0317273620702
$ws.Cells.Item(30, 15).Value = 27.8389812251787
$ws.Cells.Item(31, 15).Value = 19.0885296340919
$ws.Cells.Item(32, 15).Value = 116.66493004945
$ws.Cells.Item(33, 15).Value = 13.3178343668908
$ws.Cells.Item(34, 15).Value = 71.501738670583
$ws.Cells.Item(35, 15).Value = 16.9930181950853
$ws.Cells.Item(36, 15).Value = 14.1752090867247
$ws.Cells.Item(37, 15).Value = 14.4746086054706
$ws.Cells.Item(38, 15).Value = 14.4418266909288
$ws.Cells.Item(39, 15).Value = 25.6054366259595
$ws.Cells.Item(40, 15).Value = 10.5531774155445
$ws.Cells.Item(41, 15).Value = ".."
$ws.Cells.Item(42, 15).Value = 33.2054879022225
$ws.Cells.Item(43, 15).Value = 33.2297576485618
$ws.Cells.Item(44, 15).Value = 42.4403651619307
$ws.Cells.Item(45, 15).Value = 18.8169260722288
$ws.Cells.Item(46, 15).Value = 23.4835027649777
$ws.Cells.Item(47, 15).Value = 29.4576143998607
$ws.Cells.Item(48, 15).Value = 29.6394811757433
$ws.Cells.Item(49, 15).Value = 22.589718820485
$ws.Cells.Item(50, 15).Value = 9.13144056709746
$ws.Cells.Item(51, 15).Value = 30.027140660353
$ws.Cells.Item(52, 15).Value = 64.8059568109507
$ws.Cells.Item(53, 15).Value = 19.428377144421
$ws.Cells.Item(54, 15).Value = 32.0112993234793
$ws.Cells.Item(55, 15).Value = 27.2227142323906
$ws.Cells.Item(56, 15).Value = 15.4966004375969
$ws.Cells.Item(57, 15).Value = 11.4996686864061
$ws.Cells.Item(58, 15).Value = 23.4207724003426
$ws.Cells.Item(59, 15).Value = 26.8632747279822
$ws.Cells.Item(60, 15).Value = 22.8434868611367
$ws.Cells.Item(61, 15).Value = 17.674116448834
$ws.Cells.Item(62, 15).Value = 21.3805875037976
$ws.Cells.Item(63, 15).Value = 29.0241195092248
$ws.Cells.Item(64, 15).Value = 26.8460489152917
$ws.Cells.Item(65, 15).Value = 24.0909082505747
$ws.Cells.Item(66, 15).Value = 28.6391687409026
$ws.Cells.Item(67, 15).Value = 15.7247606941299
$ws.Cells.Item(68, 15).Value = 17.5398316378285
$ws.Cells.Item(69, 15).Value = 17.4506879069239
$ws.Cells.Item(70, 15).Value = 36.4346488607285
$ws.Cells.Item(71, 15).Value = 17.674116448834
$ws.Cells.Item(72, 15).Value = 13.1362761632726
$ws.Cells.Item(73, 15).Value = 32.3385269045631
$ws.Cells.Item(74, 15).Value = 30.5753827829633
$ws.Cells.Item(75, 15).Value = 43.3975442533836
$ws.Cells.Item(76, 15).Value = 52.9850830718795
$ws.Cells.Item(77, 15).Value = 21.2198389071567
$ws.Cells.Item(78, 15).Value = 50.2019594594326
$ws.Cells.Item(79, 15).Value = 29.9303026918056
$ws.Cells.Item(80, 15).Value = 33.4370737930998
$ws.Cells.Item(81, 15).Value = 37.650932556709
$ws.Cells.Item(82, 15).Value = 19.6853192705649
$ws.Cells.Item(83, 15).Value = 28.2433343612032
$ws.Cells.Item(84, 15).Value = 22.6728543611697
$ws.Cells.Item(85, 15).Value = 10.6337358712987
$ws.Cells.Item(86, 15).Value = 18.7646163572929
$ws.Cells.Item(87, 15).Value = 25.2099371386644
$ws.Cells.Item(88, 15).Value = 32.1852092580263
$ws.Cells.Item(89, 15).Value = 26.1971981448108
$ws.Cells.Item(90, 15).Value = 32.4220929848092
$ws.Cells.Item(91, 15).Value = 25.1953662612291
$ws.Cells.Item(92, 15).Value = 15.3101750673393
$ws.Cells.Item(93, 15).Value = 32.2817966009557
$ws.Cells.Item(94, 15).Value = 119.611789619783
$ws.Cells.Item(95, 15).Value = 20.3560475644985
$ws.Cells.Item(96, 15).Value = 31.2299145999457
$ws.Cells.Item(97, 15).Value = 18.551904407748
$ws.Cells.Item(98, 15).Value = 19.2576563826114
